# Remove the stray "Gratis" entry in B4 (the table only has A:Plato, B:Precio,
# C:Sabor columns — row 4 erroneously carried a 4th value that doesn't belong
# to any header). Clearing the cell drops it from the sheet and, since it was
# the last reference to that shared string, the shared-string table shrinks
# to match (count/uniqueCount 9 -> 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").ClearContents()
